# Implement advanced chat system with DM/Group messages, user coding,
# and participant filtering on the "Chat" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chat")

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Timestamp"
$ws.Cells.Item(1,2).Value = "Type"
$ws.Cells.Item(1,3).Value = "Participants"
$ws.Cells.Item(1,4).Value = "Sender"
$ws.Cells.Item(1,5).Value = "Message"
$ws.Cells.Item(1,6).Value = "Status"
$ws.Cells.Item(1,7).Value = "Tags"

# ---- Data rows ----
# Each entry: Timestamp, Type, Participants, Sender, Message, Status, Tags
$rows = @(
    @("20241201143000","GM","<Alyssa><Dr. Moore><Christa><Amber>","Alyssa","Hey team, where are we on the Johnson case?","active","patient-update"),
    @("20241201143100","GM","<Alyssa><Dr. Moore><Christa><Amber>","Dr. Moore","I just reviewed the medication list, all looks good","active","medical-review"),
    @("20241201143200","DM","<Alyssa><Christa>","Christa","Family meeting scheduled for tomorrow at 2pm","active","meeting"),
    @("20241201143300","DM","<Alyssa><Amber>","Alyssa","Hey Amber, can you prep the meeting notes?","active","task"),
    @("20241201143400","GM","<Alyssa><Dr. Moore><Christa><Amber>","Amber","Welcome Amber! Please connect with Alyssa on this new project","active","onboarding"),
    @("20241201143500","DM","<Dr. Moore><Christa>","Dr. Moore","Christa, can you review the Johnson medication schedule?","active","medical-task"),
    @("20241201143600","GM","<Alyssa><Dr. Moore><Christa><Amber>","Alyssa","Insurance approval came through for the Smith family!","active","good-news"),
    @("20241201143700","DM","<Alyssa><Donnie>","Alyssa","Hey Donnie, lets get Amber onboarded properly","active","onboarding"),
    @("20241201143800","GM","<Alyssa><Donnie><Amber>","Donnie","Welcome Amber! Please connect with Alyssa on this new project","active","welcome"),
    @("20241201143900","NOTE","<Alyssa>","Alyssa","Patient timeline updated - family meeting scheduled","active","patient-timeline"),
    @("20241201144000","NOTE","<Dr. Moore>","Dr. Moore","Medication review completed - no changes needed","active","medical-note"),
    @("20241201144100","DM","<Christa><Amber>","Christa","Amber, here are the key contacts for the Johnson case","active","contacts"),
    @("20241201144200","GM","<Alyssa><Dr. Moore><Christa><Amber>","Amber","Thanks everyone! Excited to be part of the team","active","introduction")
)

$r = 2
foreach ($row in $rows) {
    # Column A (Timestamp) must remain text (number stored as text), so
    # prefix with an apostrophe to prevent Excel from coercing it to a number.
    $ws.Cells.Item($r,1).Value = "'" + $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $ws.Cells.Item($r,7).Value = $row[6]
    $r = $r + 1
}
